# Update the "Product Grams " sheet: add a "Variations" column, split the
# Moringa row into per-size rows, correct the Soap / Moringa+Soap weights,
# and add the new "Moringa + Soap" grams row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Grams ")

# --- Header row ---------------------------------------------------------
# B1 used to hold "Gram"; it now holds the new "Variations" header, and the
# "Gram" header moves one column over to the new column C.
$ws.Range("B1").Value = "Variations "
$ws.Range("C1").Value = "Gram"
# Match the existing header formatting (bold, bordered, centered, top).
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4160

# New column is much wider than the rest (it holds the Gram values now).
$ws.Columns.Item(3).ColumnWidth = 66.86

# --- Data rows -----------------------------------------------------------
# Moringa used to be a single row ("100g, 200g "); it is now two rows, one
# per size, both still labelled "Moringa" in the new Variations column.
$ws.Range("A2").Value = "Moringa"
$ws.Range("B2").Value = "Moringa"
$ws.Range("C2").Value = "100g"

$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "Moringa"
$ws.Range("C3").Value = "200g"

$ws.Range("A4").Value = "Dried Curry Leaves"
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = "30g"

$ws.Range("A5").Value = "Darjeeling Black Tea "
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = "100g"

$ws.Range("A6").Value = "Combo pack"
$ws.Range("B6").ClearContents()
$ws.Range("C6").Value = "100g + 30g "

$ws.Range("A7").Value = "Soap "
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "95g"

$ws.Range("A8").Value = "Moringa + Soap "
$ws.Range("C8").Value = "100g + 95"
